$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = "2025-02-10 08:49"
$ws.Range("B11").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2018/11/07/20181107155228-740979.xls"
$ws.Range("C11").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=f9b69b43572e450dad5811f5279790ef&type=0"
$ws.Range("D11").Value = "《特种设备安全监督条例》"
$ws.Range("E11").Value = "《特种设备安全监察条例》"
$ws.Range("F11").Value = "四川省行政事业性收费标准 - 副本"
$ws.Range("G11").Value = "http://www.scnj.gov.cn/public/6598631/12263871.html"
$ws.Range("H11").Value = "四川省行政事业性收费标准"
$ws.Range("A11:H11").VerticalAlignment = -4108

# Row 12
$ws.Range("A12").Value = "2025-02-10 14:17"
$ws.Range("B12").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2020/04/07/20200407153528-264891.xlsx"
$ws.Range("C12").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=1329fefa206c474aa451fb6505e1acd4&type=0"
$ws.Range("D12").Value = "县卫生与健康局"
$ws.Range("E12").Value = "县卫生健康局"
$ws.Range("F12").Value = "南江县地方政府债务信息公开表"
$ws.Range("G12").Value = "http://www.scnj.gov.cn/public/6598531/12208091.html"
$ws.Range("H12").Value = "南江县2019年政府债务公开信息"
$ws.Range("A12:H12").VerticalAlignment = -4108

# Row 13
$ws.Range("A13").Value = "2025-02-15 10:31"
$ws.Range("B13").Value = "http://www.scnj.gov.cn/oldfiles/njx/file/p/6e4c8efc87429aa97421f0f43ba75e62.xlsx"
$ws.Range("C13").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=ef5cc5dd0c944b45b5d139f9107a2386&type=0"
$ws.Range("D13").Value = "物质储备"
$ws.Range("E13").Value = "物资储备"
$ws.Range("F13").Value = "6e4c8efc87429aa97421f0f43ba75e62.xlsx"
$ws.Range("G13").Value = "http://www.scnj.gov.cn/ywdt/gsgg/11889221.html"
$ws.Range("H13").Value = "2019年度单项目标考核情况汇总表（债务化解）"
$ws.Range("A13:H13").VerticalAlignment = -4108

# Row 14
$ws.Range("A14").Value = "2025-02-15 10:31"
$ws.Range("B14").Value = "http://www.scnj.gov.cn/oldfiles/njx/file/p/ac71be5ea5d3629cb55f831a21df8d28.xlsx"
$ws.Range("C14").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=0a2241fb656f4c77b5dcea6bf2059553&type=0"
$ws.Range("D14").Value = "物质储备"
$ws.Range("E14").Value = "物资储备"
$ws.Range("F14").Value = "ac71be5ea5d3629cb55f831a21df8d28.xlsx"
$ws.Range("G14").Value = "http://www.scnj.gov.cn/ywdt/gsgg/11889231.html"
$ws.Range("A14:G14").VerticalAlignment = -4108

# Row 15
$ws.Range("A15").Value = "2025-02-15 10:31"
$ws.Range("B15").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2020/01/14/20200114152515-753166.xlsx"
$ws.Range("C15").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=662de80e4b8245d7be097321b7fd0c0a&type=0"
$ws.Range("D15").Value = "物质储备"
$ws.Range("E15").Value = "物资储备"
$ws.Range("F15").Value = "2019年度商务局目标考核情况汇总表"
$ws.Range("G15").Value = "http://www.scnj.gov.cn/public/6598391/12142901.html"
$ws.Range("A15:G15").VerticalAlignment = -4108

# Row 16
$ws.Range("A16").Value = "2025-02-15 10:31"
$ws.Range("B16").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2020/01/10/20200110170740-461971.xls"
$ws.Range("C16").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=895d5b57fe9b416faaaed99834641e3b&type=0"
$ws.Range("D16").Value = "物质储备"
$ws.Range("E16").Value = "物资储备"
$ws.Range("F16").Value = "2019年度统计工作（含经济普查）目标考核情况公示表"
$ws.Range("G16").Value = "http://www.scnj.gov.cn/public/6598251/12100401.html"
$ws.Range("A16:G16").VerticalAlignment = -4108

# Row 17
$ws.Range("A17").Value = "2025-02-15 13:49"
$ws.Range("B17").Value = "http://www.scnj.gov.cn/oldfiles/njx/file/p/8a49594e72a63e64a83e6e25ca96b977.xlsx"
$ws.Range("C17").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=62d8fdbfda1c45bcb435ea04b2bc7a59&type=0"
$ws.Range("D17").Value = "物质储备"
$ws.Range("E17").Value = "物资储备"
$ws.Range("F17").Value = "8a49594e72a63e64a83e6e25ca96b977.xlsx"
$ws.Range("G17").Value = "http://www.scnj.gov.cn/ywdt/gsgg/11888731.html"
$ws.Range("H17").Value = "2019年度驻村帮扶考核情况汇总表"
$ws.Range("A17:H17").VerticalAlignment = -4108
